$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price-column cells whose new value is a plain decimal number
# stay stored as text (matching the sheets existing inline-string cells)
# by pre-formatting them as Text before assigning the literal string.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "34.281.70"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.791.64"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "226.63"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "0.554"
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "32.62"
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").Value = "0.294"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "2.050.76"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.799.06"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "10.99"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "0.630"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "34.304.50"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "4.25"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "68.29"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("D19").Value = "0.0₃0792"
$ws.Range("E19").Value = "  -1.70%  "
$ws.Range("D20").Value = "243.48"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "11.20"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "165.63"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").Value = "7.27"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "16.43"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "3.96"
$ws.Range("E30").Value = "  +6.21%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "3.78"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "1.400.75"
$ws.Range("E36").Value = "  -3.02%  "
$ws.Range("D37").Value = "0.666"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "84.41"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("E41").Value = "  +3.95%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "0.931"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "13.80"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "0.0523"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "1.951.09"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "104.46"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  -0.17%  "
